$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1) and "全部类型" (sheet4) both need the same set of
# "F" (想去人数) value updates, plus one "G5" cell that changes from the
# text "已售罄" to the number 218.

# Updates keyed by row number -> new F value, for the "展览" sheet.
$exhibitionUpdates = @{
    2  = 3374
    4  = 2458
    7  = 1396
    8  = 1102
    9  = 308
    13 = 102
    15 = 8636
    16 = 375
    20 = 174
    22 = 590
    24 = 1155
    26 = 2026
    29 = 1761
    31 = 1918
    35 = 91
    37 = 6
    38 = 305
    40 = 247
    41 = 429
    42 = 685
}

# Updates keyed by row number -> new F value, for the "全部类型" sheet.
$allTypesUpdates = @{
    2  = 3374
    4  = 2458
    7  = 1396
    9  = 1102
    10 = 308
    13 = 102
    15 = 8636
    16 = 375
    21 = 174
    23 = 590
    25 = 1155
    27 = 2026
    29 = 1761
    31 = 1918
    35 = 91
    37 = 6
    38 = 305
    40 = 247
    41 = 429
    46 = 685
}

$wsExhibition = $wb.Worksheets.Item("展览")
foreach ($row in $exhibitionUpdates.Keys) {
    $wsExhibition.Range("F$row").Value = $exhibitionUpdates[$row]
}
$wsExhibition.Range("G5").Value = 218

$wsAllTypes = $wb.Worksheets.Item("全部类型")
foreach ($row in $allTypesUpdates.Keys) {
    $wsAllTypes.Range("F$row").Value = $allTypesUpdates[$row]
}
$wsAllTypes.Range("G5").Value = 218
